$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new row of data for 2020-05-13 (13 Mayis 2020)
$newRow = 63
$ws.Cells.Item($newRow, 1).Value = 43964
$ws.Cells.Item($newRow, 2).Value = 33332
$ws.Cells.Item($newRow, 3).Value = 1639
$ws.Cells.Item($newRow, 4).Value = 58
$ws.Cells.Item($newRow, 5).Value = 2826

# Match style of preceding rows (date format column A, left-aligned elsewhere)
$ws.Cells.Item($newRow, 1).NumberFormat = "yyyy\-mm\-dd;@"
$ws.Cells.Item($newRow, 1).HorizontalAlignment = -4131
$ws.Cells.Item($newRow, 2).HorizontalAlignment = -4131
$ws.Cells.Item($newRow, 3).HorizontalAlignment = -4131
$ws.Cells.Item($newRow, 4).HorizontalAlignment = -4131
$ws.Cells.Item($newRow, 5).HorizontalAlignment = -4131

# Resize the table/autofilter to include the new row
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:E63"))

# Update selection to mirror the saved view state
$ws.Range("E63").Select()

